$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -------------------------------------------------
# XML 'width' attribute vs COM ColumnWidth differ by a constant offset
# (~5px/default-font = 0.8333333 chars) in this runtime; subtract it so
# the saved raw width matches the target value exactly.
$offset = 0.833333333333333

$ws.Columns.Item(4).ColumnWidth  = 12 - $offset
$ws.Columns.Item(5).ColumnWidth  = 12 - $offset
$ws.Columns.Item(8).ColumnWidth  = 20 - $offset
$ws.Columns.Item(9).ColumnWidth  = 20 - $offset
$ws.Columns.Item(10).ColumnWidth = 20 - $offset
$ws.Columns.Item(11).ColumnWidth = 20 - $offset
$ws.Columns.Item(12).ColumnWidth = 20 - $offset
$ws.Columns.Item(13).ColumnWidth = 20 - $offset
$ws.Columns.Item(14).ColumnWidth = 20 - $offset
$ws.Columns.Item(15).ColumnWidth = 20 - $offset
$ws.Columns.Item(16).ColumnWidth = 20 - $offset
$ws.Columns.Item(17).ColumnWidth = 20 - $offset
$ws.Columns.Item(18).ColumnWidth = 20 - $offset

Write-Output "columns done"

# --- Row 1: year header (2024 / 2025) ------------------------------
$ws.Range("F1:G1").UnMerge()
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G1").Value = "'2025"
$ws.Range("G1:R1").Merge()

Write-Output "row1 done"

# --- Row 2: month headers ------------------------------------------
$ws.Range("F2").Value = "'December"
$ws.Range("F3").Copy()
$ws.Range("G2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G2").Value = "'January"
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("L2").Value = "'February"
$ws.Range("P2").PasteSpecial(-4122)
$ws.Range("P2").Value = "'March"
$ws.Range("G2:K2").Merge()
$ws.Range("L2:O2").Merge()
$ws.Range("P2:R2").Merge()

Write-Output "row2 done"

# --- Row 3: week-range sub-headers ----------------------------------
$ws.Range("F3").Value = "25/Dec - 31/Dec"
$ws.Range("G3").Value = "01/Jan - 07/Jan"
$ws.Range("G3").Copy()
$ws.Range("H3:R3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H3").Value  = "08/Jan - 14/Jan"
$ws.Range("I3").Value  = "15/Jan - 21/Jan"
$ws.Range("J3").Value  = "22/Jan - 28/Jan"
$ws.Range("K3").Value  = "29/Jan - 04/Feb"
$ws.Range("L3").Value  = "05/Feb - 11/Feb"
$ws.Range("M3").Value  = "12/Feb - 18/Feb"
$ws.Range("N3").Value  = "19/Feb - 25/Feb"
$ws.Range("O3").Value  = "26/Feb - 04/Mar"
$ws.Range("P3").Value  = "05/Mar - 11/Mar"
$ws.Range("Q3").Value  = "12/Mar - 18/Mar"
$ws.Range("R3").Value  = "19/Mar - 25/Mar"

Write-Output "row3 done"

# --- Row 4: first task ("Kick Off") + new Start/End Date values ----
$ws.Range("C4").Value = "Kick Off"
$ws.Range("D4").Value = "12/25"
$ws.Range("E4").Value = "12/31"

Write-Output "row4 done"

# --- Rows 5-15: remaining tasks --------------------------------------
# Carry B4:C4's formatting down first, then fill in the per-row values.
$ws.Range("B4:C4").Copy()
$ws.Range("B5:C15").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B5").Value  = 2
$ws.Range("C5").Value  = "B"
$ws.Range("D5").Value  = "01/01"
$ws.Range("E5").Value  = "01/07"

$ws.Range("B6").Value  = 3
$ws.Range("C6").Value  = "C"
$ws.Range("D6").Value  = "01/08"
$ws.Range("E6").Value  = "01/14"

$ws.Range("B7").Value  = 4
$ws.Range("C7").Value  = "D"
$ws.Range("D7").Value  = "01/15"
$ws.Range("E7").Value  = "02/04"

$ws.Range("B8").Value  = 5
$ws.Range("C8").Value  = "E"
$ws.Range("D8").Value  = "02/05"
$ws.Range("E8").Value  = "03/04"

$ws.Range("B9").Value  = 6
$ws.Range("C9").Value  = "F"
$ws.Range("D9").Value  = "03/05"
$ws.Range("E9").Value  = "03/11"

$ws.Range("B10").Value = 7
$ws.Range("C10").Value = "Demo"
$ws.Range("D10").Value = "03/05"
$ws.Range("E10").Value = "03/11"

$ws.Range("B11").Value = 8
$ws.Range("C11").Value = "Testing"
$ws.Range("D11").Value = "03/12"
$ws.Range("E11").Value = "03/18"

$ws.Range("B12").Value = 9
$ws.Range("C12").Value = "Prod"
$ws.Range("D12").Value = "03/12"
$ws.Range("E12").Value = "03/18"

$ws.Range("B13").Value = 10
$ws.Range("D13").Value = "03/12"
$ws.Range("E13").Value = "03/18"

$ws.Range("B14").Value = 11
$ws.Range("D14").Value = "03/12"
$ws.Range("E14").Value = "03/18"

$ws.Range("B15").Value = 12
$ws.Range("D15").Value = "03/12"
$ws.Range("E15").Value = "03/18"

# Timeline highlight bars (same style as F4) -------------------------
$ws.Range("F4").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("I7:K7").PasteSpecial(-4122)
$ws.Range("L8:O8").PasteSpecial(-4122)
$ws.Range("P9").PasteSpecial(-4122)
$ws.Range("P10").PasteSpecial(-4122)
$ws.Range("Q11").PasteSpecial(-4122)
$ws.Range("Q12").PasteSpecial(-4122)
$ws.Range("Q13").PasteSpecial(-4122)
$ws.Range("Q14").PasteSpecial(-4122)
$ws.Range("Q15").PasteSpecial(-4122)

Write-Output "rows5-15 done"
